$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("optimization_parameters")
$rng = $ws.Range("A16:A17").EntireRow
$rng.Delete()
$ws.Rows.Item(16).Select()
$ws2 = $wb.Worksheets.Item("threshold_b")
$ws2.Activate()
